$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 4.469240333333333
$ws.Cells.Item(2, 8).Value = 13.407721
$ws.Cells.Item(2, 9).Value = 0.3358648218165975
$ws.Cells.Item(2, 10).Value = 0.3358648218165975
$ws.Cells.Item(2, 13).Value = 3.795192333333334
$ws.Cells.Item(2, 14).Value = 11.385577
$ws.Cells.Item(2, 15).Value = 0.01044213755712683
$ws.Cells.Item(2, 16).Value = 0.01044213755712683
$ws.Cells.Item(2, 17).Value = 16.96162664889078
$ws.Cells.Item(2, 18).Value = 152.654639840017
$ws.Cells.Item(2, 19).Value = 0.003507146670008805
$ws.Cells.Item(2, 20).Value = 0.003507146670008805

$ws.Cells.Item(3, 7).Value = 4.469240333333333
$ws.Cells.Item(3, 8).Value = 13.407721
$ws.Cells.Item(3, 9).Value = 0.3358648218165975
$ws.Cells.Item(3, 10).Value = 0.3358648218165975
$ws.Cells.Item(3, 14).Value = 730.1291960000001
$ws.Cells.Item(3, 15).Value = 0.6696287328350964
$ws.Cells.Item(3, 16).Value = 0.6696287328350964
$ws.Cells.Item(3, 17).Value = 1087.707617102479
$ws.Cells.Item(3, 18).Value = 9789.368553922317
$ws.Cells.Item(3, 19).Value = 0.2249047350369336
$ws.Cells.Item(3, 20).Value = 0.2249047350369336

$ws.Cells.Item(4, 7).Value = 4.469240333333333
$ws.Cells.Item(4, 8).Value = 13.407721
$ws.Cells.Item(4, 9).Value = 0.3358648218165975
$ws.Cells.Item(4, 10).Value = 0.3358648218165975
$ws.Cells.Item(4, 13).Value = 29.801371
$ws.Cells.Item(4, 14).Value = 89.404113
$ws.Cells.Item(4, 15).Value = 0.08199584844219236
$ws.Cells.Item(4, 16).Value = 0.08199584844219235
$ws.Cells.Item(4, 17).Value = 133.1894892618303
$ws.Cells.Item(4, 18).Value = 1198.705403356473
$ws.Cells.Item(4, 19).Value = 0.02753952102673768
$ws.Cells.Item(4, 20).Value = 0.02753952102673767

$ws.Cells.Item(5, 7).Value = 4.469240333333333
$ws.Cells.Item(5, 8).Value = 13.407721
$ws.Cells.Item(5, 9).Value = 0.3358648218165975
$ws.Cells.Item(5, 10).Value = 0.3358648218165975
$ws.Cells.Item(5, 13).Value = 86.47679266666667
$ws.Cells.Item(5, 14).Value = 259.430378
$ws.Cells.Item(5, 15).Value = 0.2379332811655844
$ws.Cells.Item(5, 16).Value = 0.2379332811655844
$ws.Cells.Item(5, 17).Value = 386.4855696831708
$ws.Cells.Item(5, 18).Value = 3478.370127148538
$ws.Cells.Item(5, 19).Value = 0.07991341908291741
$ws.Cells.Item(5, 20).Value = 0.07991341908291741

$ws.Cells.Item(6, 9).Value = 0.2721973992379558
$ws.Cells.Item(6, 10).Value = 0.2721973992379558
$ws.Cells.Item(6, 13).Value = 3.795192333333334
$ws.Cells.Item(6, 14).Value = 11.385577
$ws.Cells.Item(6, 15).Value = 0.01044213755712683
$ws.Cells.Item(6, 16).Value = 0.01044213755712683
$ws.Cells.Item(6, 17).Value = 13.74633590889845
$ws.Cells.Item(6, 18).Value = 123.717023180086
$ws.Cells.Item(6, 19).Value = 0.002842322685534905
$ws.Cells.Item(6, 20).Value = 0.002842322685534905

$ws.Cells.Item(7, 9).Value = 0.2721973992379558
$ws.Cells.Item(7, 10).Value = 0.2721973992379558
$ws.Cells.Item(7, 14).Value = 730.1291960000001
$ws.Cells.Item(7, 15).Value = 0.6696287328350964
$ws.Cells.Item(7, 16).Value = 0.6696287328350964
$ws.Cells.Item(7, 17).Value = 881.518888775681
$ws.Cells.Item(7, 18).Value = 7933.66999898113
$ws.Cells.Item(7, 19).Value = 0.1822711995327212
$ws.Cells.Item(7, 20).Value = 0.1822711995327212

$ws.Cells.Item(8, 9).Value = 0.2721973992379558
$ws.Cells.Item(8, 10).Value = 0.2721973992379558
$ws.Cells.Item(8, 13).Value = 29.801371
$ws.Cells.Item(8, 14).Value = 89.404113
$ws.Cells.Item(8, 15).Value = 0.08199584844219236
$ws.Cells.Item(8, 16).Value = 0.08199584844219235
$ws.Cells.Item(8, 17).Value = 107.9417379492593
$ws.Cells.Item(8, 18).Value = 971.4756415433339
$ws.Cells.Item(8, 19).Value = 0.02231905669427435
$ws.Cells.Item(8, 20).Value = 0.02231905669427434

$ws.Cells.Item(9, 9).Value = 0.2721973992379558
$ws.Cells.Item(9, 10).Value = 0.2721973992379558
$ws.Cells.Item(9, 13).Value = 86.47679266666667
$ws.Cells.Item(9, 14).Value = 259.430378
$ws.Cells.Item(9, 15).Value = 0.2379332811655844
$ws.Cells.Item(9, 16).Value = 0.2379332811655844
$ws.Cells.Item(9, 17).Value = 313.2223444591783
$ws.Cells.Item(9, 18).Value = 2819.001100132604
$ws.Cells.Item(9, 19).Value = 0.06476482032542537
$ws.Cells.Item(9, 20).Value = 0.06476482032542535

$ws.Cells.Item(10, 7).Value = 3.887787333333333
$ws.Cells.Item(10, 8).Value = 11.663362
$ws.Cells.Item(10, 9).Value = 0.2921684453243378
$ws.Cells.Item(10, 10).Value = 0.2921684453243377
$ws.Cells.Item(10, 13).Value = 3.795192333333334
$ws.Cells.Item(10, 14).Value = 11.385577
$ws.Cells.Item(10, 15).Value = 0.01044213755712683
$ws.Cells.Item(10, 16).Value = 0.01044213755712683
$ws.Cells.Item(10, 17).Value = 14.75490068109711
$ws.Cells.Item(10, 18).Value = 132.794106129874
$ws.Cells.Item(10, 19).Value = 0.003050863095928625
$ws.Cells.Item(10, 20).Value = 0.003050863095928624

$ws.Cells.Item(11, 7).Value = 3.887787333333333
$ws.Cells.Item(11, 8).Value = 11.663362
$ws.Cells.Item(11, 9).Value = 0.2921684453243378
$ws.Cells.Item(11, 10).Value = 0.2921684453243377
$ws.Cells.Item(11, 14).Value = 730.1291960000001
$ws.Cells.Item(11, 15).Value = 0.6696287328350964
$ws.Cells.Item(11, 16).Value = 0.6696287328350964
$ws.Cells.Item(11, 17).Value = 946.1956799685503
$ws.Cells.Item(11, 18).Value = 8515.761119716952
$ws.Cells.Item(11, 19).Value = 0.1956443858169364
$ws.Cells.Item(11, 20).Value = 0.1956443858169364

$ws.Cells.Item(12, 7).Value = 3.887787333333333
$ws.Cells.Item(12, 8).Value = 11.663362
$ws.Cells.Item(12, 9).Value = 0.2921684453243378
$ws.Cells.Item(12, 10).Value = 0.2921684453243377
$ws.Cells.Item(12, 13).Value = 29.801371
$ws.Cells.Item(12, 14).Value = 89.404113
$ws.Cells.Item(12, 15).Value = 0.08199584844219236
$ws.Cells.Item(12, 16).Value = 0.08199584844219235
$ws.Cells.Item(12, 17).Value = 115.8613926897673
$ws.Cells.Item(12, 18).Value = 1042.752534207906
$ws.Cells.Item(12, 19).Value = 0.02395659956240537
$ws.Cells.Item(12, 20).Value = 0.02395659956240536

$ws.Cells.Item(13, 7).Value = 3.887787333333333
$ws.Cells.Item(13, 8).Value = 11.663362
$ws.Cells.Item(13, 9).Value = 0.2921684453243378
$ws.Cells.Item(13, 10).Value = 0.2921684453243377
$ws.Cells.Item(13, 13).Value = 86.47679266666667
$ws.Cells.Item(13, 14).Value = 259.430378
$ws.Cells.Item(13, 15).Value = 0.2379332811655844
$ws.Cells.Item(13, 16).Value = 0.2379332811655844
$ws.Cells.Item(13, 17).Value = 336.2033791567596
$ws.Cells.Item(13, 18).Value = 3025.830412410836
$ws.Cells.Item(13, 19).Value = 0.06951659684906733
$ws.Cells.Item(13, 20).Value = 0.06951659684906732

$ws.Cells.Item(14, 7).Value = 1.327597
$ws.Cells.Item(14, 8).Value = 3.982791
$ws.Cells.Item(14, 9).Value = 0.09976933362110896
$ws.Cells.Item(14, 10).Value = 0.09976933362110893
$ws.Cells.Item(14, 13).Value = 3.795192333333334
$ws.Cells.Item(14, 14).Value = 11.385577
$ws.Cells.Item(14, 15).Value = 0.01044213755712683
$ws.Cells.Item(14, 16).Value = 0.01044213755712683
$ws.Cells.Item(14, 17).Value = 5.038485956156334
$ws.Cells.Item(14, 18).Value = 45.346373605407
$ws.Cells.Item(14, 19).Value = 0.001041805105654499
$ws.Cells.Item(14, 20).Value = 0.001041805105654498

$ws.Cells.Item(15, 7).Value = 1.327597
$ws.Cells.Item(15, 8).Value = 3.982791
$ws.Cells.Item(15, 9).Value = 0.09976933362110896
$ws.Cells.Item(15, 10).Value = 0.09976933362110893
$ws.Cells.Item(15, 14).Value = 730.1291960000001
$ws.Cells.Item(15, 15).Value = 0.6696287328350964
$ws.Cells.Item(15, 16).Value = 0.6696287328350964
$ws.Cells.Item(15, 17).Value = 323.1057767406706
$ws.Cells.Item(15, 18).Value = 2907.951990666036
$ws.Cells.Item(15, 19).Value = 0.06680841244850518
$ws.Cells.Item(15, 20).Value = 0.06680841244850515

$ws.Cells.Item(16, 7).Value = 1.327597
$ws.Cells.Item(16, 8).Value = 3.982791
$ws.Cells.Item(16, 9).Value = 0.09976933362110896
$ws.Cells.Item(16, 10).Value = 0.09976933362110893
$ws.Cells.Item(16, 13).Value = 29.801371
$ws.Cells.Item(16, 14).Value = 89.404113
$ws.Cells.Item(16, 15).Value = 0.08199584844219236
$ws.Cells.Item(16, 16).Value = 0.08199584844219235
$ws.Cells.Item(16, 17).Value = 39.56421073548699
$ws.Cells.Item(16, 18).Value = 356.077896619383
$ws.Cells.Item(16, 19).Value = 0.008180671158774976
$ws.Cells.Item(16, 20).Value = 0.008180671158774973

$ws.Cells.Item(17, 7).Value = 1.327597
$ws.Cells.Item(17, 8).Value = 3.982791
$ws.Cells.Item(17, 9).Value = 0.09976933362110896
$ws.Cells.Item(17, 10).Value = 0.09976933362110893
$ws.Cells.Item(17, 13).Value = 86.47679266666667
$ws.Cells.Item(17, 14).Value = 259.430378
$ws.Cells.Item(17, 15).Value = 0.2379332811655844
$ws.Cells.Item(17, 16).Value = 0.2379332811655844
$ws.Cells.Item(17, 17).Value = 114.8063305138887
$ws.Cells.Item(17, 18).Value = 1033.256974624998
$ws.Cells.Item(17, 19).Value = 0.02373844490817431
$ws.Cells.Item(17, 20).Value = 0.02373844490817431
